# Auto update Excel log
#
# Appends freshly-logged sensor rows produced on 2026-02-06:
#   - "ALERTS" sheet gains rows 4-5 (a MODERATE then a CRITICAL bathroom
#     no-motion alert).
#   - "PIR" sheet gains rows 71-77 (raw PIR motion/no-motion readings for
#     the bathroom sensor).

function Add-LogRow {
    param($ws, $row, $date, $timestamp, $hour, $location, $value, $status)

    # Column A holds a literal "YYYY-MM-DD" string, which Excel would
    # otherwise auto-convert into a date serial number. Force the cell to
    # Text first so the value round-trips as plain text...
    $dateCell = $ws.Cells.Item($row, 1)
    $dateCell.NumberFormat = "@"
    $dateCell.Value = $date
    # ...then drop the cell back to the workbook's default style so no
    # stray number-format style lingers on the cell (matches the rest of
    # the log, which carries no explicit style).
    $dateCell.Style = "Normal"

    $ws.Cells.Item($row, 2).Value = $timestamp
    $ws.Cells.Item($row, 3).Value = $hour
    $ws.Cells.Item($row, 4).Value = $location
    $ws.Cells.Item($row, 5).Value = $value
    $ws.Cells.Item($row, 6).Value = $status
}

$wb = $excel.ActiveWorkbook

# --- ALERTS sheet: append rows 4-5 ---
$alerts = $wb.Worksheets.Item("ALERTS")

Add-LogRow $alerts 4 "2026-02-06" "09:42:50" "09:00" "Bathroom" "MODERATE" "MODERATE ALERT: Bathroom occupied, no motion > 40s."
Add-LogRow $alerts 5 "2026-02-06" "09:43:10" "09:00" "Bathroom" "CRITICAL" "CRITICAL ALERT: Bathroom occupied, no motion > 60s."

# --- PIR sheet: append rows 71-77 ---
$pir = $wb.Worksheets.Item("PIR")

Add-LogRow $pir 71 "2026-02-06" "09:43:06" "09:00" "Bathroom" "No Motion" "Inactive"
Add-LogRow $pir 72 "2026-02-06" "09:43:11" "09:00" "Bathroom" "No Motion" "Inactive"
Add-LogRow $pir 73 "2026-02-06" "09:43:16" "09:00" "Bathroom" "No Motion" "Inactive"
Add-LogRow $pir 74 "2026-02-06" "09:43:21" "09:00" "Bathroom" "No Motion" "Inactive"
Add-LogRow $pir 75 "2026-02-06" "09:43:26" "09:00" "Bathroom" "No Motion" "Inactive"
Add-LogRow $pir 76 "2026-02-06" "09:43:29" "09:00" "Bathroom" "Motion Detected" "Active"
Add-LogRow $pir 77 "2026-02-06" "09:43:37" "09:00" "Bathroom" "No Motion" "Inactive"
